# JM added ProII files
$wb = $excel.ActiveWorkbook

# --- Update view state on the "A+" sheet first (it will stop being the active tab) ---
$wsA = $wb.Worksheets.Item("A+")
$wsA.Range("A1:D1").Select() | Out-Null

# --- Add the new "ProII" worksheet after the last existing sheet ("A+") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ProII"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 17.42578125
$ws.Columns.Item(3).ColumnWidth = 13.5703125
$ws.Columns.Item(4).ColumnWidth = 33.140625

# --- Header row ---
$ws.Range("A1").Value = "TestName"
$ws.Range("B1").Value = "Class"
$ws.Range("C1").Value = "OID"
$ws.Range("D1").Value = "Note"

# --- Data rows ---
$ws.Range("A2").Value = "Simple HX"
$ws.Range("B2").Value = "HeatExchanger"
$ws.Range("C2").Value = 316

$ws.Range("A3").Value = "Pump"
$ws.Range("B3").Value = "Pump"
$ws.Range("C3").Value = 345

$ws.Range("A4").Value = "Valve"
$ws.Range("B4").Value = "Valve"
$ws.Range("C4").Value = 299

$ws.Range("A5").Value = "Flash"
$ws.Range("B5").Value = "Separator"
$ws.Range("C5").Value = 309

# --- Selection / view state for the new sheet (becomes the active tab) ---
$ws.Range("C15").Select() | Out-Null
